# Cucumber/Jenkins reporting run results appended to the "Output" sheet.
# Rows 2-11 get refreshed with the results of a newer test run (timestamps
# and figures from 23/01/2022), and the old wrap-text styling that was
# applied to the "row does not exist" messages (rows 8-11, column C) is
# cleared along with the oversized row-8 height.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output")

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric/comma-separated strings (e.g. "97,900")
    # aren't silently coerced into numbers by the COM value setter.
    $rng.NumberFormat = "@"
    $rng.Value = $value
    # Re-align the style with an already-unstyled neighbour cell so the
    # cell doesn't end up pinned to a custom "Text" number-format style.
    $rng.Style = $ws.Range("A1").Style
}

# Row 2
$ws.Range("A2").Value = "23/01/2022 2:23:53 pm"
$ws.Range("D2").Value = "₹1,22,140"
$ws.Range("E2").Value = "₹8,368"
$ws.Range("G2").Value = "₹1,30,518"

# Row 3
$ws.Range("A3").Value = "23/01/2022 2:24:00 pm"

# Row 4 (gains Base Fare/Fee/Addons/Total columns this run)
$ws.Range("A4").Value = "23/01/2022 2:24:42 pm"
$ws.Range("B4").Value = "N"
$ws.Range("C4").Value = "N/A"
Set-TextValue "D4" "1,032,360"
Set-TextValue "E4" "41,082"
Set-TextValue "F4" "10"
Set-TextValue "G4" "1,073,452"

# Row 5
$ws.Range("A5").Value = "23/01/2022 2:24:51 pm"

# Row 6
$ws.Range("A6").Value = "23/01/2022 2:25:03 pm"
$ws.Range("D6").Value = "₹1,45,040"
$ws.Range("E6").Value = "₹10,360"
$ws.Range("G6").Value = "₹1,55,410"

# Row 7 (gains Base Fare/Fee/Addons/Total columns this run)
$ws.Range("A7").Value = "23/01/2022 2:25:18 pm"
$ws.Range("B7").Value = "N"
$ws.Range("C7").Value = "N/A"
Set-TextValue "D7" "97,900"
Set-TextValue "E7" "8,980"
Set-TextValue "F7" "10"
Set-TextValue "G7" "106,890"

# Row 8 - also drop the wrap-text style + oversized row height from the
# earlier "Autosuggest element..." long message that used to live here.
$ws.Range("A8").Value = "23/01/2022 2:25:33 pm"
$ws.Range("C8").Value = "Row 6 is empty"
$ws.Range("C8").Style = $ws.Range("B8").Style
$ws.Rows(8).AutoFit()

# Row 9 - drop wrap-text style
$ws.Range("A9").Value = "23/01/2022 2:25:38 pm"
$ws.Range("C9").Value = "Row 7 does not exist"
$ws.Range("C9").Style = $ws.Range("B9").Style

# Row 10 - drop wrap-text style
$ws.Range("A10").Value = "23/01/2022 2:25:43 pm"
$ws.Range("C10").Value = "Row 8 does not exist"
$ws.Range("C10").Style = $ws.Range("B10").Style

# Row 11 - drop wrap-text style
$ws.Range("A11").Value = "23/01/2022 2:25:48 pm"
$ws.Range("C11").Value = "Row 9 does not exist"
$ws.Range("C11").Style = $ws.Range("B11").Style
